$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update E13: Quickbook/Xero intro text (now with "as" and theme bold run)
# ---------------------------------------------------------------------------
$ws.Range("E13").Value = "It should show 'Integrate with other online services , Accounting as Quickbooks Online and Xero'"
$ws.Range("E13").Characters(16, 81).Font.Bold = $true

# ---------------------------------------------------------------------------
# 2. Insert two new rows (14 and 15) for the Quickbooks / Xero testcases,
#    pushing the old rows 14-15 down to 16-17.
# ---------------------------------------------------------------------------
$ws.Rows("14:15").Insert()

# Row heights / row-level formatting to match s="8" customFormat rows (like rows 10-11)
$ws.Rows(14).RowHeight = 78.75
$ws.Rows(15).RowHeight = 78.75

# Pull the cell formatting (font/border/alignment) for columns A-D from row 10,
# which already uses the same visual style (s=8/2/2/7) that the new rows need.
$ws.Range("A10:D11").Copy() | Out-Null
$ws.Range("A14:D15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Row 14 : Quickbooks settings test ---
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "Suppliers Home page"
$ws.Cells.Item(14, 3).Value = "View Company settings page"
$ws.Cells.Item(14, 4).Value = "Company settings->Integration->Quickbooks"

$ws.Range("E14").Value = 'It shows "Integration settings as Outlet, Export as Select Outlet dropdown, Category, Tax code Tracking class as Please select dropdown and Cancel and Save"'
$ws.Range("E14").Characters(9, 148).Font.Bold = $true

$ws.Cells.Item(14, 6).Value = "It gets displayed the Quickbooks updated OK"

# --- Row 15 : Xero settings test ---
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "Suppliers Home page"
$ws.Cells.Item(15, 3).Value = "View Company settings page"
$ws.Cells.Item(15, 4).Value = "Company settings->Integration->Xero"

$ws.Range("E15").Value = 'It shows "Outlet, Export as Select Customer dropdown, Category as Please select dropdown and Cancel and Save"'
$ws.Range("E15").Characters(9, 101).Font.Bold = $true

$ws.Range("F15").Value = "It gets displayed the Xero  OK"
$ws.Range("F15").Characters(22, 9).Font.Bold = $true

Write-Host "Rows 14/15 populated"

# ---------------------------------------------------------------------------
# 2b. Fix the "SL. No" numbering on the rows that were shifted down
#     (old row 14 -> row 16 [was 13, now 15], old row 15 -> row 17 [was 14, now 16])
# ---------------------------------------------------------------------------
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(17, 1).Value = 16

# ---------------------------------------------------------------------------
# 3. Fix up the selection / active cell to the bottom of the sheet
# ---------------------------------------------------------------------------
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A17").Select() | Out-Null

Write-Host "Done"
